# Update "想去人数" (F column) counts across the 4 sheets of the workbook
# per the data refresh recorded in the commit diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 269
$ws1.Range("F5").Value = 329
$ws1.Range("F7").Value = 2189
$ws1.Range("F10").Value = 1637
$ws1.Range("F11").Value = 1637
$ws1.Range("F13").Value = 64
$ws1.Range("F14").Value = 1413
$ws1.Range("F17").Value = 582
$ws1.Range("F20").Value = 7263
$ws1.Range("F21").Value = 8006
$ws1.Range("F22").Value = 49
$ws1.Range("F27").Value = 93
$ws1.Range("F30").Value = 255
$ws1.Range("F35").Value = 1444
$ws1.Range("F36").Value = 208
$ws1.Range("F41").Value = 736
$ws1.Range("F43").Value = 1365
$ws1.Range("F46").Value = 200
$ws1.Range("F47").Value = 87

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 60
$ws2.Range("F17").Value = 8

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2624

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 269
$ws4.Range("F7").Value = 329
$ws4.Range("F10").Value = 2189
$ws4.Range("F13").Value = 1637
$ws4.Range("F14").Value = 1637
$ws4.Range("F15").Value = 64
$ws4.Range("F16").Value = 1413
$ws4.Range("F18").Value = 582
$ws4.Range("F21").Value = 60
$ws4.Range("F24").Value = 7263
$ws4.Range("F25").Value = 8006
$ws4.Range("F26").Value = 49
$ws4.Range("F28").Value = 93
$ws4.Range("F31").Value = 1444
$ws4.Range("F32").Value = 208
$ws4.Range("F39").Value = 736
$ws4.Range("F43").Value = 1365
$ws4.Range("F46").Value = 200
$ws4.Range("F48").Value = 8
